# Update ASV_rank (column H) values to reflect refreshed comparison counts.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = "316/422"
$ws.Range("H12").Value = "312/422"
$ws.Range("H15").Value = "343/422"
$ws.Range("H17").Value = "289/422"
$ws.Range("H20").Value = "269/422"
$ws.Range("H23").Value = "138/422"
$ws.Range("H24").Value = "369/422"
$ws.Range("H28").Value = "388/422"
$ws.Range("H33").Value = "178/422"
$ws.Range("H38").Value = "304/422"
$ws.Range("H41").Value = "393/422"
$ws.Range("H42").Value = "225/422"
$ws.Range("H43").Value = "359/422"
$ws.Range("H45").Value = "310/422"
$ws.Range("H47").Value = "411/422"
$ws.Range("H50").Value = "303/422"
$ws.Range("H56").Value = "296/422"
$ws.Range("H62").Value = "384/422"
$ws.Range("H71").Value = "222/422"
$ws.Range("H73").Value = "342/422"
$ws.Range("H76").Value = "282/422"
$ws.Range("H82").Value = "137/422"
$ws.Range("H85").Value = "248/422"
$ws.Range("H87").Value = "386/422"
$ws.Range("H89").Value = "337/422"
$ws.Range("H90").Value = "364/422"
$ws.Range("H91").Value = "366/422"
$ws.Range("H92").Value = "367/422"
$ws.Range("H93").Value = "380/422"
